$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2765957446808511
$ws.Range("C2").Value = 0.3829787234042553
$ws.Range("J2").Value = 0.02127659574468085
$ws.Range("P2").Value = 0.2127659574468085
$ws.Range("S2").Value = 0.1063829787234043
$ws.Range("J3").Value = 0.05882352941176471
$ws.Range("P3").Value = 0.5294117647058824
$ws.Range("S3").Value = 0.4117647058823529
$ws.Range("E6").Value = 0.03125
$ws.Range("F6").Value = 0.0625
$ws.Range("J6").Value = 0.1875
$ws.Range("Q6").Value = 0.125
$ws.Range("R6").Value = 0.125
$ws.Range("S6").Value = 0.46875
$ws.Range("B7").Value = 0.1282051282051282
$ws.Range("F7").Value = 0.02564102564102564
$ws.Range("J7").Value = 0.2051282051282051
$ws.Range("Q7").Value = 0.07692307692307693
$ws.Range("R7").Value = 0.02564102564102564
$ws.Range("S7").Value = 0.5384615384615384
$ws.Range("B8").Value = 0.08181818181818182
$ws.Range("F8").Value = 0.05454545454545454
$ws.Range("J8").Value = 0.1090909090909091
$ws.Range("O8").Value = 0.01818181818181818
$ws.Range("Q8").Value = 0.1
$ws.Range("R8").Value = 0.08181818181818182
$ws.Range("S8").Value = 0.5545454545454546
$ws.Range("B9").Value = 0.03225806451612903
$ws.Range("D9").Value = 0.03225806451612903
$ws.Range("F9").Value = 0.06451612903225806
$ws.Range("J9").Value = 0.03225806451612903
$ws.Range("Q9").Value = 0.06451612903225806
$ws.Range("S9").Value = 0.7741935483870968
$ws.Range("B10").Value = 0.1137724550898204
$ws.Range("F10").Value = 0.07784431137724551
$ws.Range("J10").Value = 0.125748502994012
$ws.Range("O10").Value = 0.01796407185628742
$ws.Range("Q10").Value = 0.08982035928143713
$ws.Range("R10").Value = 0.08383233532934131
$ws.Range("S10").Value = 0.4910179640718563
$ws.Range("G11").Value = 0.1639344262295082
$ws.Range("J11").Value = 0.1147540983606557
$ws.Range("K11").Value = 0.2459016393442623
$ws.Range("L11").Value = 0.4098360655737705
$ws.Range("S11").Value = 0.06557377049180328
$ws.Range("G12").Value = 0.88
$ws.Range("J12").Value = 0.12
$ws.Range("G13").Value = 0.75
$ws.Range("J13").Value = 0.25
$ws.Range("H15").Value = 0.1363636363636364
$ws.Range("I15").Value = 0.1363636363636364
$ws.Range("J15").Value = 0.2272727272727273
$ws.Range("K15").Value = 0.1363636363636364
$ws.Range("M15").Value = 0.04545454545454546
$ws.Range("O15").Value = 0.04545454545454546
$ws.Range("S15").Value = 0.2727272727272727
$ws.Range("H16").Value = 0.3157894736842105
$ws.Range("J16").Value = 0.3684210526315789
$ws.Range("K16").Value = 0.1578947368421053
$ws.Range("O16").Value = 0.05263157894736842
$ws.Range("S16").Value = 0.1052631578947368
$ws.Range("F17").Value = 0.02857142857142857
$ws.Range("H17").Value = 0.1714285714285714
$ws.Range("I17").Value = 0.05714285714285714
$ws.Range("J17").Value = 0.4
$ws.Range("K17").Value = 0.1142857142857143
$ws.Range("S17").Value = 0.2285714285714286
$ws.Range("H18").Value = 0.1428571428571428
$ws.Range("I18").Value = 0.07142857142857142
$ws.Range("J18").Value = 0.4642857142857143
$ws.Range("K18").Value = 0.1785714285714286
$ws.Range("M18").Value = 0.03571428571428571
$ws.Range("S18").Value = 0.1071428571428571
$ws.Range("F19").Value = 0.01481481481481482
$ws.Range("H19").Value = 0.337037037037037
$ws.Range("I19").Value = 0.08888888888888889
$ws.Range("J19").Value = 0.262962962962963
$ws.Range("K19").Value = 0.1111111111111111
$ws.Range("M19").Value = 0.02222222222222222
$ws.Range("N19").Value = 0.003703703703703704
$ws.Range("O19").Value = 0.03703703703703703
$ws.Range("S19").Value = 0.1222222222222222
